# ACE Score Calculator - "Conditional Formatting - Yellow"
#
# 1. B16 (months/years-of-education selector) changes from 13 to 9 on both
#    "ACE Scoring Calculator" and "ACE scoring_Gen" sheets. All of the
#    dependent formulas in C17:H25 recalculate automatically.
# 2. The B17:B25 conditional formatting on "ACE Scoring Calculator" gains a
#    third rule ("Minor Impairment" -> yellow/gold) alongside the existing
#    "No Impairments" (green) and "Major impairment" (red) rules.
# 3. Selection / active-cell bookkeeping is refreshed on all three sheets.

$wb = $excel.ActiveWorkbook

$wsCalc = $wb.Worksheets.Item("ACE Scoring Calculator")
$wsGen  = $wb.Worksheets.Item("ACE scoring_Gen")
$wsAce  = $wb.Worksheets.Item("ACE scoring")

# ---------------------------------------------------------------------------
# 1. Update the education-years selector -> cascades through C17:H25.
# ---------------------------------------------------------------------------
$wsCalc.Range("B16").Value = 9
$wsGen.Range("B16").Value = 9

# ---------------------------------------------------------------------------
# 2. Rebuild the conditional formatting on "ACE Scoring Calculator"!B17:B25
#    so "Minor Impairment" (yellow/gold) is represented too, alongside the
#    pre-existing "No Impairments" (green) and "Major impairment" (red)
#    rules. The three formulas test mutually exclusive states of H17, so
#    relative rule priority has no visible effect; ordering mirrors the
#    sheet's existing No Impairments / Major impairment / Minor Impairment
#    presentation order.
# ---------------------------------------------------------------------------
$cfRange = $wsCalc.Range("B17:B25")
$cfRange.FormatConditions.Delete()

# COM colors are packed as BGR (0xBBGGRR), not RGB.
$colorGreen = 0x50 * 65536 + 0xD0 * 256 + 0x92   # 92D050 green        -> No Impairments
$colorRed   = 0x00 * 65536 + 0x00 * 256 + 0xFF   # FF0000 red          -> Major impairment
$colorGold  = 0x00 * 65536 + 0xC0 * 256 + 0xFF   # FFC000 gold/yellow  -> Minor Impairment

$ruleNone = $cfRange.FormatConditions.Add(2, 0, '=$H17="No Impairments"')
$ruleNone.Interior.Color = $colorGreen

$ruleMajor = $cfRange.FormatConditions.Add(2, 0, '=$H17="Major impairment"')
$ruleMajor.Interior.Color = $colorRed

$ruleMinor = $cfRange.FormatConditions.Add(2, 0, '=$H17="Minor Impairment"')
$ruleMinor.Interior.Color = $colorGold

# ---------------------------------------------------------------------------
# 3. Refresh the saved selection on every sheet.
# ---------------------------------------------------------------------------
$wsAce.Activate()
$wsAce.Range("F35").Select()

$wsGen.Activate()
$wsGen.Range("B17").Select()

$wsCalc.Activate()
$wsCalc.Range("H29").Select()
